$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 34,9
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0.01358695652173913
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0.01380175658720201
$arr[0,8] = 0.02171471359041559
$arr[1,0] = 0
$arr[1,1] = 0.001872659176029963
$arr[1,2] = 0.02173913043478261
$arr[1,3] = 0
$arr[1,4] = 0
$arr[1,5] = 0.1483488132094948
$arr[1,6] = 0
$arr[1,7] = 0.02383939774153074
$arr[1,8] = 0.03594159490827404
$arr[2,0] = 0
$arr[2,1] = 0.03121098626716608
$arr[2,2] = 0.01358695652173913
$arr[2,3] = 0
$arr[2,4] = 0.01351351351351351
$arr[2,5] = 0.005159958720330238
$arr[2,6] = 0.0066711140760507
$arr[2,7] = 0.007528230865746549
$arr[2,8] = 0.005990265818045677
$arr[3,0] = 1
$arr[3,1] = 0
$arr[3,2] = 0.358695652173914
$arr[3,3] = 0.6666666666666666
$arr[3,4] = 0.3918918918918917
$arr[3,5] = 0.001031991744066047
$arr[3,6] = 0
$arr[3,7] = 0.2383939774153082
$arr[3,8] = 0.01010857356795208
$arr[4,0] = 0
$arr[4,1] = 0
$arr[4,2] = 0.01358695652173913
$arr[4,3] = 0
$arr[4,4] = 0
$arr[4,5] = 0.01238390092879257
$arr[4,6] = 0
$arr[4,7] = 0.0150564617314931
$arr[4,8] = 0.01647323099962562
$arr[5,0] = 0
$arr[5,1] = 0.02933832709113611
$arr[5,2] = 0
$arr[5,3] = 0
$arr[5,4] = 0
$arr[5,5] = 0.003611971104231167
$arr[5,6] = 0.04536357571714479
$arr[5,7] = 0
$arr[5,8] = 0.0007487832272557095
$arr[6,0] = 0
$arr[6,1] = 0.1186017478152312
$arr[6,2] = 0.06793478260869562
$arr[6,3] = 0
$arr[6,4] = 0.04054054054054054
$arr[6,5] = 0.03070175438596482
$arr[6,6] = 0.1467645096731154
$arr[6,7] = 0.06273525721455465
$arr[6,8] = 0.0378135529764133
$arr[7,0] = 0
$arr[7,1] = 0.07303370786516868
$arr[7,2] = 0
$arr[7,3] = 0
$arr[7,4] = 0
$arr[7,5] = 0.007739938080495358
$arr[7,6] = 0.09673115410273539
$arr[7,7] = 0.001254705144291092
$arr[7,8] = 0.002620741295394983
$arr[8,0] = 0
$arr[8,1] = 0.01373283395755306
$arr[8,2] = 0
$arr[8,3] = 0
$arr[8,4] = 0.01351351351351351
$arr[8,5] = 0.001289989680082559
$arr[8,6] = 0.07671781187458321
$arr[8,7] = 0.005018820577164366
$arr[8,8] = 0.001871958068139274
$arr[9,0] = 0
$arr[9,1] = 0.003121098626716604
$arr[9,2] = 0
$arr[9,3] = 0
$arr[9,4] = 0
$arr[9,5] = 0
$arr[9,6] = 0.01400933955970646
$arr[9,7] = 0
$arr[9,8] = 0
$arr[10,0] = 0
$arr[10,1] = 0.01435705368289638
$arr[10,2] = 0
$arr[10,3] = 0
$arr[10,4] = 0
$arr[10,5] = 0.0005159958720330237
$arr[10,6] = 0.02735156771180785
$arr[10,7] = 0
$arr[10,8] = 0.0003743916136278548
$arr[11,0] = 0
$arr[11,1] = 0
$arr[11,2] = 0.03260869565217391
$arr[11,3] = 0.3333333333333333
$arr[11,4] = 0.06756756756756757
$arr[11,5] = 0
$arr[11,6] = 0
$arr[11,7] = 0.03136762860727729
$arr[11,8] = 0.006364657431673532
$arr[12,0] = 0
$arr[12,1] = 0
$arr[12,2] = 0.02445652173913043
$arr[12,3] = 0
$arr[12,4] = 0
$arr[12,5] = 0
$arr[12,6] = 0
$arr[12,7] = 0.09535759096612294
$arr[12,8] = 0.001123174840883564
$arr[13,0] = 0
$arr[13,1] = 0
$arr[13,2] = 0
$arr[13,3] = 0
$arr[13,4] = 0
$arr[13,5] = 0
$arr[13,6] = 0
$arr[13,7] = 0
$arr[13,8] = 0
$arr[14,0] = 0
$arr[14,1] = 0.1004993757802749
$arr[14,2] = 0.07608695652173909
$arr[14,3] = 0
$arr[14,4] = 0.05405405405405406
$arr[14,5] = 0.04643962848297194
$arr[14,6] = 0.04936624416277523
$arr[14,7] = 0.08531994981179424
$arr[14,8] = 0.08199176338450015
$arr[15,0] = 0
$arr[15,1] = 0
$arr[15,2] = 0
$arr[15,3] = 0
$arr[15,4] = 0
$arr[15,5] = 0.01160990712074304
$arr[15,6] = 0
$arr[15,7] = 0.003764115432873275
$arr[15,8] = 0.01946836390864846
$arr[16,0] = 0
$arr[16,1] = 0.0006242197253433209
$arr[16,2] = 0
$arr[16,3] = 0
$arr[16,4] = 0
$arr[16,5] = 0
$arr[16,6] = 0.01134089392928619
$arr[16,7] = 0
$arr[16,8] = 0
$arr[17,0] = 0
$arr[17,1] = 0
$arr[17,2] = 0.09510869565217386
$arr[17,3] = 0
$arr[17,4] = 0.1216216216216216
$arr[17,5] = 0
$arr[17,6] = 0
$arr[17,7] = 0.03387703889585948
$arr[17,8] = 0
$arr[18,0] = 0
$arr[18,1] = 0.01061173533083645
$arr[18,2] = 0.01358695652173913
$arr[18,3] = 0
$arr[18,4] = 0
$arr[18,5] = 0.0327657378740969
$arr[18,6] = 0.00200133422281521
$arr[18,7] = 0.02634880803011292
$arr[18,8] = 0.1078247847248224
$arr[19,0] = 0
$arr[19,1] = 0.08177278401997519
$arr[19,2] = 0.03532608695652173
$arr[19,3] = 0
$arr[19,4] = 0.02702702702702703
$arr[19,5] = 0.04385964912280684
$arr[19,6] = 0.01667778519012674
$arr[19,7] = 0.03262233375156838
$arr[19,8] = 0.05428678397603877
$arr[20,0] = 0
$arr[20,1] = 0
$arr[20,2] = 0
$arr[20,3] = 0
$arr[20,4] = 0
$arr[20,5] = 0
$arr[20,6] = 0
$arr[20,7] = 0
$arr[20,8] = 0
$arr[21,0] = 0
$arr[21,1] = 0.0961298377028716
$arr[21,2] = 0.02989130434782608
$arr[21,3] = 0
$arr[21,4] = 0
$arr[21,5] = 0.124355005159958
$arr[21,6] = 0.0153435623749166
$arr[21,7] = 0.05520702634880809
$arr[21,8] = 0.05540995881692232
$arr[22,0] = 0
$arr[22,1] = 0.002496878901373283
$arr[22,2] = 0
$arr[22,3] = 0
$arr[22,4] = 0
$arr[22,5] = 0
$arr[22,6] = 0.02001334222815209
$arr[22,7] = 0
$arr[22,8] = 0
$arr[23,0] = 0
$arr[23,1] = 0.05305867665418236
$arr[23,2] = 0
$arr[23,3] = 0
$arr[23,4] = 0
$arr[23,5] = 0.007481940144478846
$arr[23,6] = 0.00266844563042028
$arr[23,7] = 0
$arr[23,8] = 0.002620741295394983
$arr[24,0] = 0
$arr[24,1] = 0
$arr[24,2] = 0
$arr[24,3] = 0
$arr[24,4] = 0
$arr[24,5] = 0
$arr[24,6] = 0.00066711140760507
$arr[24,7] = 0
$arr[24,8] = 0
$arr[25,0] = 0
$arr[25,1] = 0
$arr[25,2] = 0
$arr[25,3] = 0
$arr[25,4] = 0
$arr[25,5] = 0
$arr[25,6] = 0.03602401601067377
$arr[25,7] = 0
$arr[25,8] = 0
$arr[26,0] = 0
$arr[26,1] = 0.03682896379525598
$arr[26,2] = 0
$arr[26,3] = 0
$arr[26,4] = 0
$arr[26,5] = 0.0128998968008256
$arr[26,6] = 0
$arr[26,7] = 0
$arr[26,8] = 0.001497566454511419
$arr[27,0] = 0
$arr[27,1] = 0.01435705368289638
$arr[27,2] = 0
$arr[27,3] = 0
$arr[27,4] = 0
$arr[27,5] = 0.001805985552115583
$arr[27,6] = 0.02001334222815209
$arr[27,7] = 0
$arr[27,8] = 0.001497566454511419
$arr[28,0] = 0
$arr[28,1] = 0.009987515605493132
$arr[28,2] = 0
$arr[28,3] = 0
$arr[28,4] = 0
$arr[28,5] = 0.002837977296181631
$arr[28,6] = 0.04869913275517015
$arr[28,7] = 0
$arr[28,8] = 0.002995132909022838
$arr[29,0] = 0
$arr[29,1] = 0.02372034956304621
$arr[29,2] = 0
$arr[29,3] = 0
$arr[29,4] = 0
$arr[29,5] = 0.002837977296181631
$arr[29,6] = 0.0426951300867245
$arr[29,7] = 0
$arr[29,8] = 0
$arr[30,0] = 0
$arr[30,1] = 0.1598002496878905
$arr[30,2] = 0.01630434782608696
$arr[30,3] = 0
$arr[30,4] = 0.04054054054054054
$arr[30,5] = 0.03663570691434456
$arr[30,6] = 0.0720480320213477
$arr[30,7] = 0.02509410288582183
$arr[30,8] = 0.02283788843129915
$arr[31,0] = 0
$arr[31,1] = 0.05305867665418236
$arr[31,2] = 0
$arr[31,3] = 0
$arr[31,4] = 0
$arr[31,5] = 0.008255933952528381
$arr[31,6] = 0.1667778519012671
$arr[31,7] = 0.003764115432873275
$arr[31,8] = 0.005241482590789967
$arr[32,0] = 0
$arr[32,1] = 0
$arr[32,2] = 0.002717391304347826
$arr[32,3] = 0
$arr[32,4] = 0
$arr[32,5] = 0
$arr[32,6] = 0
$arr[32,7] = 0.0150564617314931
$arr[32,8] = 0.01123174840883565
$arr[33,0] = 0
$arr[33,1] = 0.0006242197253433209
$arr[33,2] = 0
$arr[33,3] = 0
$arr[33,4] = 0
$arr[33,5] = 0
$arr[33,6] = 0.04069379586390928
$arr[33,7] = 0
$arr[33,8] = 0

$ws.Range("B2:J35").Value = $arr
